$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 28
$ws.Range("A28").Value = "20251121-005"
$ws.Range("D28").Value = "QUO-20251121-003"
$ws.Range("E28").Value = 580

# Row 29
$ws.Range("D29").Value = "QUO-20251121-005"
$ws.Range("E29").Value = 1015

# Row 30
$ws.Range("A30").Value = "20251121-007"
$ws.Range("D30").Value = "QUO-20251121-008"

# Row 31
$ws.Range("D31").Value = "QUO-20251121-001"
$ws.Range("E31").Value = 7540

# Row 32
$ws.Range("D32").Value = "QUO-20251121-004"
$ws.Range("E32").Value = 7540
